$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name on both sheets (append "-1st")
$wsInput.Range("B1").Value  = "2525-MS-EPP-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-INST-1st"
$wsOutput.Range("B1").Value = "2525-MS-EPP-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-1-DATE-VAR-INST-1st"

# shortname / description changed from numeric 2525 to text "252d"
$wsInput.Range("B2").Value = "252d"
$wsInput.Range("B3").Value = "252d"

# Reset selection on the input sheet back to the top
$wsInput.Range("B1").Select()

# Make the output sheet the active / selected tab
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
